# Auto-generated Excel COM-interop script
# Updates market-data-derived profit columns (H-N) on the 8 "Gungnir_Profits" sheets
# per the scheduled runner refresh (chore: update Sheets via scheduled runner).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 58.8125
$ws.Range("I8").Value = 58.8125
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 176.4375
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -37.4375
$ws.Range("N8").ClearContents()

$ws.Range("H64").Value = 3143.75
$ws.Range("I64").Value = 2810
$ws.Range("J64").Value = 3700
$ws.Range("K64").Value = 2810
$ws.Range("L64").Value = 3700
$ws.Range("M64").Value = -2562
$ws.Range("N64").Value = -4196

$ws.Range("H67").Value = 3143.75
$ws.Range("I67").Value = 2810
$ws.Range("J67").Value = 3700
$ws.Range("K67").Value = 2810
$ws.Range("L67").Value = 3700
$ws.Range("M67").Value = -1952
$ws.Range("N67").Value = -5416

$ws.Range("H95").Value = 30625
$ws.Range("J95").Value = 30625
$ws.Range("L95").Value = 30625
$ws.Range("N95").Value = -36117

$ws.Range("H129").Value = 917.67926
$ws.Range("J129").Value = 1027.1428
$ws.Range("L129").Value = 3081.4284
$ws.Range("N129").Value = -13081.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 15000
$ws.Range("J6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("N6").Value = -15346

$ws.Range("H32").Value = 12348675
$ws.Range("I32").Value = 1294799.6
$ws.Range("J32").Value = 250007000
$ws.Range("K32").Value = 1294799.6
$ws.Range("L32").Value = 250007000
$ws.Range("M32").Value = -1294512.6
$ws.Range("N32").Value = -250007574

$ws.Range("H62").Value = 17479.5
$ws.Range("J62").Value = 17479.5
$ws.Range("L62").Value = 17479.5
$ws.Range("N62").Value = -18727.5

$ws.Range("H65").Value = 17479.5
$ws.Range("J65").Value = 17479.5
$ws.Range("L65").Value = 52438.5
$ws.Range("N65").Value = -58678.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 20542
$ws.Range("J28").Value = 20542
$ws.Range("L28").Value = 20542
$ws.Range("N28").Value = -21130

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H64").Value = 39800
$ws.Range("J64").Value = 39800
$ws.Range("L64").Value = 39800
$ws.Range("N64").Value = -40296

$ws.Range("H67").Value = 39800
$ws.Range("J67").Value = 39800
$ws.Range("L67").Value = 39800
$ws.Range("N67").Value = -41516

$ws.Range("H95").Value = 8314.4
$ws.Range("J95").Value = 8314.4
$ws.Range("L95").Value = 8314.4
$ws.Range("N95").Value = -13806.4

$ws.Range("H107").Value = 727.6667
$ws.Range("I107").Value = 585.53845
$ws.Range("J107").Value = 1097.2
$ws.Range("K107").Value = 585.53845
$ws.Range("L107").Value = 1097.2
$ws.Range("M107").Value = 1334.46155
$ws.Range("N107").Value = -4937.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 150.06667
$ws.Range("I6").Value = 112.583336
$ws.Range("K6").Value = 337.750008
$ws.Range("M6").Value = -224.750008

$ws.Range("H10").Value = 560
$ws.Range("I10").Value = 150
$ws.Range("J10").Value = 833.3333
$ws.Range("K10").Value = 450
$ws.Range("L10").Value = 2499.9999
$ws.Range("M10").Value = -311
$ws.Range("N10").Value = -2777.9999

$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2251
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -5256
$ws.Range("N66").ClearContents()

$ws.Range("H87").Value = 7980
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 7980
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 23940
$ws.Range("N87").Value = -26436
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 7980
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 7980
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 71820
$ws.Range("N90").Value = -84300
$ws.Range("M90").ClearContents()

$ws.Range("H131").Value = 862.0204
$ws.Range("J131").Value = 862.3505
$ws.Range("L131").Value = 2587.0515
$ws.Range("N131").Value = -12667.0515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 24130.5
$ws.Range("J39").Value = 24130.5
$ws.Range("L39").Value = 24130.5
$ws.Range("N39").Value = -25194.5

$ws.Range("H52").Value = 22999
$ws.Range("I52").Value = 6000
$ws.Range("J52").Value = 28665.334
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 28665.334
$ws.Range("M52").Value = -5741
$ws.Range("N52").Value = -29183.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 8483.666999999999
$ws.Range("I62").Value = 8483.666999999999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 8483.666999999999
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -7859.666999999999
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 12766.667
$ws.Range("J64").Value = 12766.667
$ws.Range("L64").Value = 12766.667
$ws.Range("N64").Value = -13216.667

$ws.Range("H65").Value = 8483.666999999999
$ws.Range("I65").Value = 8483.666999999999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25451.001
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -22331.001
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 12766.667
$ws.Range("J67").Value = 12766.667
$ws.Range("L67").Value = 12766.667
$ws.Range("N67").Value = -14326.667

$ws.Range("H74").Value = 21695
$ws.Range("I74").Value = 25500
$ws.Range("J74").Value = 17890
$ws.Range("K74").Value = 25500
$ws.Range("L74").Value = 17890
$ws.Range("M74").Value = -24502
$ws.Range("N74").Value = -19886

$ws.Range("H77").Value = 21695
$ws.Range("I77").Value = 25500
$ws.Range("J77").Value = 17890
$ws.Range("K77").Value = 76500
$ws.Range("L77").Value = 53670
$ws.Range("M77").Value = -71508
$ws.Range("N77").Value = -63654

$ws.Range("H132").Value = 21282924
$ws.Range("I132").Value = 55557572
$ws.Range("J132").Value = 9003.379000000001
$ws.Range("K132").Value = 166672716
$ws.Range("L132").Value = 27010.137
$ws.Range("M132").Value = -166670186
$ws.Range("N132").Value = -32070.137

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2936.3635
$ws.Range("I62").Value = 2787.5
$ws.Range("K62").Value = 2787.5
$ws.Range("M62").Value = -2163.5

$ws.Range("H63").Value = 19249
$ws.Range("J63").Value = 19249
$ws.Range("L63").Value = 19249
$ws.Range("N63").Value = -20497

$ws.Range("H65").Value = 2936.3635
$ws.Range("I65").Value = 2787.5
$ws.Range("K65").Value = 13937.5
$ws.Range("M65").Value = -10817.5

$ws.Range("H66").Value = 19249
$ws.Range("J66").Value = 19249
$ws.Range("L66").Value = 57747
$ws.Range("N66").Value = -63987

$ws.Range("H80").Value = 14500
$ws.Range("J80").Value = 14500
$ws.Range("L80").Value = 14500
$ws.Range("N80").Value = -16496

$ws.Range("H83").Value = 14500
$ws.Range("J83").Value = 14500
$ws.Range("L83").Value = 43500
$ws.Range("N83").Value = -53484

$ws.Range("H92").Value = 49650
$ws.Range("J92").Value = 49650
$ws.Range("L92").Value = 49650
$ws.Range("N92").Value = -54642

$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492

$ws.Range("H97").Value = 25000
$ws.Range("J97").Value = 25000
$ws.Range("L97").Value = 25000
$ws.Range("N97").Value = -26982

$ws.Range("H109").Value = 20341.8
$ws.Range("J109").Value = 20341.8
$ws.Range("L109").Value = 20341.8
$ws.Range("N109").Value = -23115.8
